# Generate Report for Handoff
# Re-runs the localization status report generation: swaps the old e2e test
# file's guid for a new one (and its content hash), refreshes the "latest
# handoff" timestamps, and clears out the per-language "already handed
# back" bookkeeping (Latest Target File / Latest Handback File / Latest
# Handback DateTime) now that this is a fresh handoff instead of a
# finished handback.

$wb = $excel.ActiveWorkbook

$oldGuid = "637ea68b-9530-496a-8634-572befa58fe0"
$newGuid = "0ca6ae1a-240e-475a-837d-852025b0b7e3"
$oldHash = "b4dccdb3cbd5d4f2873307003f8f6c4628faa669"
$newHash = "c1abc839fc5f466652ffa3c04fbda87a64baad26"

$oldFileName   = "$oldGuid.md"
$newFileName   = "$newGuid.md"
$oldPath       = "e2e\$oldGuid.md"
$newPath       = "e2e\$newGuid.md"
$oldZhXlf      = "$oldGuid.$oldHash.zh-cn.xlf"
$newZhXlf      = "$newGuid.$newHash.zh-cn.xlf"
$oldDeXlf      = "$oldGuid.$oldHash.de-de.xlf"
$newDeXlf      = "$newGuid.$newHash.de-de.xlf"

$noHandback = "0001-01-01 00:00:00"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPath
$wsOverview.Range("G2").Value = "2016-09-04 21:06:56"

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$2') {
        $hl.TextToDisplay = $newPath
    }
}

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newFileName
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = "2016-09-04 21:06:51"
$wsZh.Range("K2").Value = $noHandback

foreach ($hl in $wsZh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newFileName
    } elseif ($addr -eq '$I$2') {
        $hl.Delete()
    }
}

$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("I2").ClearFormats()
$wsZh.Columns.Item(9).AutoFit()
$wsZh.Columns.Item(10).AutoFit()

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newFileName
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = "2016-09-04 21:06:51"
$wsDe.Range("K2").Value = $noHandback

foreach ($hl in $wsDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newFileName
    } elseif ($addr -eq '$I$2') {
        $hl.Delete()
    }
}

$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("I2").ClearFormats()
$wsDe.Columns.Item(9).AutoFit()
$wsDe.Columns.Item(10).AutoFit()
